# The commit reverts the deck's "Integral" theme color scheme back to the
# default Office Theme color scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# Helper: PowerPoint/VBA's RGB() packs r,g,b into an OLE_COLOR integer
# (0x00BBGGRR) - replicate that here since this host doesn't expose RGB().
function ToOle([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# The presentation's theme color scheme (shared by the slide master / the
# whole deck) - edit its 12 slots in place so it matches the stock
# "Office Theme" color scheme.
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = ToOle 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = ToOle 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = ToOle 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = ToOle 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = ToOle 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = ToOle 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = ToOle 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = ToOle 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = ToOle 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = ToOle 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = ToOle 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = ToOle 0x95 0x4F 0x72   # folHlink 954F72
